$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range (A1:B3) before writing the new, smaller table
$ws.Range("A1:B3").ClearContents()

# Write the new single-row task set
$ws.Range("A1").Value = "2 2 5"
$ws.Range("B1").Value = "3 5 7"
$ws.Range("C1").Value = "1 5 10"

# Update the active selection to match the new last-used cell
$ws.Range("C1").Select()
